$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text values (e.g. "1.003", "23.950.10").
# Excel auto-converts plain numeric-looking strings assigned via .Value into
# real numbers, so force a temporary text format on the specific cells being
# updated, assign the new text, then restore each cell's original (default)
# number format by pasting formats back in from an untouched reference cell.
$dCells = @("D2","D3","D5","D6","D7","D8","D9","D10","D11","D12","D13","D14","D15","D16","D17","D18","D19","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($addr in $dCells) {
  $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "23.950.10"
$ws.Range("D3").Value = "1.647.92"
$ws.Range("D5").Value = "310.12"
$ws.Range("D6").Value = "1.002"
$ws.Range("D7").Value = "0.3886"
$ws.Range("D8").Value = "0.3797"
$ws.Range("D9").Value = "51.90"
$ws.Range("D10").Value = "1.343"
$ws.Range("D11").Value = "1.003"
$ws.Range("D12").Value = "0.08435"
$ws.Range("D13").Value = "23.83"
$ws.Range("D14").Value = "7.039"
$ws.Range("D15").Value = "7.990"
$ws.Range("D16").Value = "0.00001307"
$ws.Range("D17").Value = "1.650.08"
$ws.Range("D18").Value = "94.26"
$ws.Range("D19").Value = "0.06986"
$ws.Range("D20").Value = "19.63"
$ws.Range("D21").Value = "6.949"
$ws.Range("D22").Value = "1.002"
$ws.Range("D23").Value = "13.73"
$ws.Range("D24").Value = "23.951.34"
$ws.Range("D25").Value = "2.449"
$ws.Range("D26").Value = "2.928"
$ws.Range("D27").Value = "21.99"
$ws.Range("D28").Value = "152.65"
$ws.Range("D29").Value = "5.398"
$ws.Range("D30").Value = "137.79"
$ws.Range("D31").Value = "7.812"
$ws.Range("D32").Value = "2.521"
$ws.Range("D33").Value = "1.830.26"
$ws.Range("D34").Value = "1.017"
$ws.Range("D35").Value = "0.08032"
$ws.Range("D36").Value = "6.710"
$ws.Range("D37").Value = "0.02907"
$ws.Range("D38").Value = "0.2669"
$ws.Range("D39").Value = "10.70"
$ws.Range("D40").Value = "0.09072"
$ws.Range("D41").Value = "0.7571"
$ws.Range("D42").Value = "13.36"
$ws.Range("D43").Value = "1.418"
$ws.Range("D44").Value = "16.08"
$ws.Range("D45").Value = "0.6949"
$ws.Range("D46").Value = "2.443"
$ws.Range("D47").Value = "4.091"
$ws.Range("D48").Value = "1.002"
$ws.Range("D49").Value = "0.08305"
$ws.Range("D50").Value = "134.74"
$ws.Range("D51").Value = "1.225"

$ws.Range("E2").Value = "  -2.12%  "
$ws.Range("E3").Value = "  -1.18%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("E5").Value = "  -1.21%  "
$ws.Range("E6").Value = "  +0.29%  "
$ws.Range("E7").Value = "  -2.07%  "
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  -4.73%  "
$ws.Range("E12").Value = "  -1.94%  "
$ws.Range("E13").Value = "  -2.50%  "
$ws.Range("E14").Value = "  -4.29%  "
$ws.Range("E15").Value = "  +0.93%  "
$ws.Range("E16").Value = "  -3.22%  "
$ws.Range("E17").Value = "  -1.05%  "
$ws.Range("E18").Value = "  -1.31%  "
$ws.Range("E19").Value = "  +0.01%  "
$ws.Range("E20").Value = "  -5.01%  "
$ws.Range("E21").Value = "  -1.10%  "
$ws.Range("E22").Value = "  +0.33%  "
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("E24").Value = "  -2.11%  "
$ws.Range("E25").Value = "  +0.84%  "
$ws.Range("E26").Value = "  -4.30%  "
$ws.Range("E27").Value = "  -2.48%  "
$ws.Range("E28").Value = "  -3.26%  "
$ws.Range("E29").Value = "  -1.12%  "
$ws.Range("E30").Value = "  -3.84%  "
$ws.Range("E31").Value = "  -4.83%  "
$ws.Range("E32").Value = "  -0.89%  "
$ws.Range("E33").Value = "  -1.01%  "
$ws.Range("E34").Value = "  -4.38%  "
$ws.Range("E35").Value = "  -3.25%  "
$ws.Range("E36").Value = "  -2.79%  "
$ws.Range("E37").Value = "  -4.57%  "
$ws.Range("E38").Value = "  -4.01%  "
$ws.Range("E39").Value = "  -4.97%  "
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("E41").Value = "  -2.57%  "
$ws.Range("E42").Value = "  -4.07%  "
$ws.Range("E43").Value = "  -2.01%  "
$ws.Range("E44").Value = "  -3.23%  "
$ws.Range("E45").Value = "  -2.90%  "
$ws.Range("E46").Value = "  -4.16%  "
$ws.Range("E47").Value = "  -1.37%  "
$ws.Range("E48").Value = "  +0.26%  "
$ws.Range("E49").Value = "  -1.89%  "
$ws.Range("E50").Value = "  -1.50%  "
$ws.Range("E51").Value = "  -4.27%  "

# Restore the default (General) cell format/style on the updated D cells so
# only the cell contents change, matching the source edit.
$ws.Range("F1").Copy()
foreach ($addr in $dCells) {
  $ws.Range($addr).PasteSpecial(-4122)
}
$ws.Range("F1").ClearContents()
$excel.CutCopyMode = $false

